# research_proposal_biblio.xlsx edit script
# "correction climate change + link with body mass as a LHT"
#
# Fills in bibliography rows 71-74 with four new references (climate-change /
# elevation-warming papers, plus a life-history-theory / body-mass paper),
# copying the F-column conditional formatting (the "~"/"no" styling) from
# existing rows so the new cells pick up the same cell style indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 71 : Ohmura (2012) ----------------------------------------------
$ws.Range("A71").Value = "Enhanced temperature variability in high-altitude climate change"
$ws.Range("B71").Value = "Theoritical and Applied Climatology"
$ws.Range("C71").Value = 2012
$ws.Range("D71").Value = "Ohmura"
$ws.Range("F71").Value = "≈"
$ws.Range("G71").Value = "Temperature variability in high altitude"
$ws.Range("I71").Value = "yes"
$ws.Range("J71").Value = "yes"
$ws.Range("K71").Value = "yes"

# ---- Row 72 : Kittel, Thornton, Royle, Chase (2002) -----------------------
$ws.Range("A72").Value = "Climates of the Rocky Mountains: Historical and Future Patterns"
$ws.Range("B72").Value = "Ø"
$ws.Range("C72").Value = 2002
$ws.Range("D72").Value = "Kittel, Thornton, Royle, Chase"
$ws.Range("E72").Value = "Book chapter"
$ws.Range("F72").Value = "≈"
$ws.Range("G72").Value = "Droughts"
$ws.Range("I72").Value = "yes"
$ws.Range("J72").Value = "yes"
$ws.Range("K72").Value = "yes"

# ---- Row 73 : Giorgi, Hurrell, Marinucci, Beniston (1997) -----------------
$ws.Range("A73").Value = "Elevation Dependency of the Surface Climate Change Signal: A Model Study"
$ws.Range("B73").Value = "Journal of Climate"
$ws.Range("C73").Value = 1997
$ws.Range("D73").Value = "Giorgi, Hurrell, Marinucci, Beniston"
$ws.Range("F73").Value = "≈"
$ws.Range("G73").Value = """more pronounced warming at high elevation"""
$ws.Range("I73").Value = "yes"
$ws.Range("J73").Value = "yes"
$ws.Range("K73").Value = "yes"

# ---- Row 74 : Bell (1980) --------------------------------------------------
$ws.Range("A74").Value = "The Costs of Reproduction and Their Consequences"
$ws.Range("B74").Value = "The american naturalist"
$ws.Range("C74").Value = 1980
$ws.Range("D74").Value = "Bell"
$ws.Range("F74").Value = "no"
$ws.Range("G74").Value = "LHT"
$ws.Range("I74").Value = "yes"
$ws.Range("J74").Value = "yes"
$ws.Range("K74").Value = "yes"

# Copy the existing F-column conditional formatting onto the new cells so the
# "~" / "no" values keep the same look used throughout the rest of the sheet.
$ws.Range("F69").Copy() | Out-Null
$ws.Range("F71:F73").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F19").Copy() | Out-Null
$ws.Range("F74").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = $false

# ---- View state: scroll position / active selection -----------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 62
$win.ScrollColumn = 5
$ws.Range("K74").Select() | Out-Null
